$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from "Students" to "Sheet1"
$ws.Name = "Sheet1"

# Column A ("Id") held numeric ids as text ("2","4","5"); store as real numbers
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 5

# Column E ("Subscribed") held "True"/"False" as text; store as real booleans
$ws.Range("E2").Value = $true
$ws.Range("E3").Value = $true
$ws.Range("E4").Value = $false
